$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the Offense text in C2, then fill down for the remaining
# applicable records (rows 3-10) - "copied down for applicable records (1 to many)"
$ws.Range("C2").Value = "Arrest - Illegal Weapon Possesions"
$ws.Range("C2:C10").FillDown() | Out-Null

# Auto-fit the new Offense column to the entered text
$ws.Columns.Item(3).AutoFit() | Out-Null

# Leave the selection where the author left off
$ws.Range("D7").Select() | Out-Null
